$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 280; D = 44995; J = 240; K = 3500; L = 3500; M = 3500; O = "Región Metropolitana"; P = 1167 },
    @{ Row = 281; D = 44161; J = 80; K = 3000; L = 3000; M = 3000; O = "Región Metropolitana"; P = 1000 },
    @{ Row = 282; D = 44721; J = 80; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 283; D = 44655; J = 40; K = 6000; L = 6000; M = 6000; O = "Región Metropolitana"; P = 2000 },
    @{ Row = 284; D = 44351; J = 240; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 285; D = 44553; J = 120; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 286; D = 44974; J = 240; K = 3500; L = 3500; M = 3500; O = "Región Metropolitana"; P = 1167 },
    @{ Row = 287; D = 44365; J = 240; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 288; D = 44529; J = 80; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 289; D = 44323; J = 240; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 290; D = 44963; J = 40; K = 3500; L = 3500; M = 3500; O = "Región Metropolitana"; P = 1167 },
    @{ Row = 291; D = 44966; J = 80; K = 3500; L = 3500; M = 3500; O = "Región Metropolitana"; P = 1167 },
    @{ Row = 292; D = 44495; J = 240; K = 2500; L = 3000; M = 2750; O = "Región Metropolitana"; P = 917 },
    @{ Row = 293; D = 44987; J = 80; K = 3500; L = 3500; M = 3500; O = "Región Metropolitana"; P = 1167 },
    @{ Row = 294; D = 44763; J = 80; K = 4500; L = 4500; M = 4500; O = "Región Metropolitana"; P = 1500 },
    @{ Row = 295; D = 44574; J = 120; K = 3000; L = 3000; M = 3000; O = "Región Metropolitana"; P = 1000 },
    @{ Row = 296; D = 44383; J = 240; K = 3000; L = 3000; M = 3000; O = "Región Metropolitana"; P = 1000 },
    @{ Row = 297; D = 44701; J = 240; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 298; D = 44516; J = 240; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 299; D = 44806; J = 240; K = 3000; L = 3500; M = 3250; O = "Región Metropolitana"; P = 1083 },
    @{ Row = 300; D = 44880; J = 240; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 301; D = 44530; J = 320; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 302; D = 44511; J = 80; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 303; D = 44242; J = 50; K = 5000; L = 5000; M = 5000; O = "Provincia de Cautín"; P = 1667 },
    @{ Row = 304; D = 44637; J = 80; K = 4000; L = 4000; M = 4000; O = "Región Metropolitana"; P = 1333 },
    @{ Row = 305; D = 44771; J = 240; K = 4000; L = 4000; M = 4000; O = "Región Metropolitana"; P = 1333 },
    @{ Row = 306; D = 44847; J = 80; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 307; D = 44665; J = 80; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 308; D = 44371; J = 80; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 309; D = 44882; J = 120; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 310; D = 44257; J = 200; K = 6000; L = 6000; M = 6000; O = "Provincia de Cautín"; P = 2000 },
    @{ Row = 311; D = 44278; J = 120; K = 5000; L = 5000; M = 5000; O = "Provincia de Cautín"; P = 1667 },
    @{ Row = 312; D = 44544; J = 280; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 313; D = 44432; J = 240; K = 5000; L = 5000; M = 5000; O = "Región Metropolitana"; P = 1667 },
    @{ Row = 314; D = 44784; J = 80; K = 4000; L = 4000; M = 4000; O = "Región Metropolitana"; P = 1333 },
    @{ Row = 315; D = 44750; J = 240; K = 3500; L = 3500; M = 3500; O = "Región Metropolitana"; P = 1167 },
    @{ Row = 316; D = 44483; J = 120; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 317; D = 44945; J = 80; K = 3000; L = 3000; M = 3000; O = "Región Metropolitana"; P = 1000 },
    @{ Row = 318; D = 44747; J = 240; K = 2500; L = 3000; M = 2750; O = "Región Metropolitana"; P = 917 },
    @{ Row = 319; D = 44610; J = 240; K = 2500; L = 3000; M = 2750; O = "Región Metropolitana"; P = 917 },
    @{ Row = 320; D = 44819; J = 160; K = 3000; L = 3000; M = 3000; O = "Región Metropolitana"; P = 1000 },
    @{ Row = 321; D = 44951; J = 40; K = 6000; L = 6000; M = 6000; O = "Provincia de Cautín"; P = 2000 },
    @{ Row = 322; D = 44812; J = 120; K = 3000; L = 3000; M = 3000; O = "Región Metropolitana"; P = 1000 },
    @{ Row = 323; D = 44876; J = 240; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 324; D = 44370; J = 40; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 325; D = 44447; J = 40; K = 4500; L = 4500; M = 4500; O = "Región Metropolitana"; P = 1500 },
    @{ Row = 326; D = 44859; J = 240; K = 2000; L = 2500; M = 2250; O = "Región Metropolitana"; P = 750 },
    @{ Row = 327; D = 44467; J = 240; K = 3500; L = 3500; M = 3500; O = "Región Metropolitana"; P = 1167 },
    @{ Row = 328; D = 44532; J = 120; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana"; P = 833 },
    @{ Row = 329; D = 44799; J = 240; K = 4000; L = 4000; M = 4000; O = "Región Metropolitana"; P = 1333 },
    @{ Row = 330; D = 44327; J = 260; K = 2500; L = 3000; M = 2769; O = "Región Metropolitana"; P = 923 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value2 = $u.D
    $ws.Cells.Item($u.Row, 10).Value2 = $u.J
    $ws.Cells.Item($u.Row, 11).Value2 = $u.K
    $ws.Cells.Item($u.Row, 12).Value2 = $u.L
    $ws.Cells.Item($u.Row, 13).Value2 = $u.M
    $ws.Cells.Item($u.Row, 15).Value2 = $u.O
    $ws.Cells.Item($u.Row, 16).Value2 = $u.P
}

# Row 330 is brand new - fill in the columns that stay constant across the block,
# copied from row 329 (which itself is unaffected in columns A,B,C,E,F,G,H,I,N,Q,R).
$ws.Cells.Item(330, 1).Value2 = $ws.Cells.Item(329, 1).Value2
$ws.Cells.Item(330, 2).Value2 = $ws.Cells.Item(329, 2).Value2
$ws.Cells.Item(330, 3).Value2 = $ws.Cells.Item(329, 3).Value2
$ws.Cells.Item(330, 5).Value2 = $ws.Cells.Item(329, 5).Value2
$ws.Cells.Item(330, 6).Value2 = $ws.Cells.Item(329, 6).Value2
$ws.Cells.Item(330, 7).Value2 = $ws.Cells.Item(329, 7).Value2
$ws.Cells.Item(330, 8).Value2 = $ws.Cells.Item(329, 8).Value2
$ws.Cells.Item(330, 9).Value2 = $ws.Cells.Item(329, 9).Value2
$ws.Cells.Item(330, 14).Value2 = $ws.Cells.Item(329, 14).Value2
$ws.Cells.Item(330, 17).Value2 = $ws.Cells.Item(329, 17).Value2
$ws.Cells.Item(330, 18).Value2 = $ws.Cells.Item(329, 18).Value2

# Copy the date number format from D329 to the newly created D330 cell.
$ws.Cells.Item(330, 4).NumberFormat = $ws.Cells.Item(329, 4).NumberFormat
